$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash the existing header style (A1, style index 1 - the highlighted fill)
# in a scratch cell so we can re-apply it after the old data is cleared.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats

# Wipe the old LOGINCREDS table (values + formatting)
$ws.Range("A1:C4").Clear()

# New header row + single data row
$ws.Range("B1").Value = "mobilnumber"
$ws.Range("C1").Value = "partnerloanid"
$ws.Range("D1").Value = "pasrtnercustomerid"
$ws.Range("B2").Value = 88833

# Re-apply the stashed header style onto the new header cells
$ws.Range("F1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$ws.Range("F1").Clear()

# Column widths
$ws.Columns.Item(1).ColumnWidth = 4
$ws.Columns.Item(2).ColumnWidth = 24
$ws.Columns.Item(3).ColumnWidth = 23.5
$ws.Columns.Item(4).ColumnWidth = 17.833333333333332

# Selection / active cell
$ws.Range("B2").Select()
